$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "groups" sheet
$ws2 = $wb.Worksheets.Item(2)   # "group-student" sheet

# ---------------------------------------------------------------------------
# 1. Restructure "groups" sheet: swap the GROUP_NAME / ASSIGNMENT_NAME
#    columns (C<->D) using Cut+Insert so the exact (bestFit) column widths
#    travel with the data, then add two new trailing columns for the
#    per-student info that used to live only on the "group-student" sheet.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(4).Cut()
$ws1.Columns.Item(3).Insert()

$ws1.Range("E1").Value = "STUDENT_ID"
$ws1.Range("F1").Value = "IS_MANAGER"

# ---------------------------------------------------------------------------
# 2. Fill in one row per student (merging in the group-student roster),
#    replacing the old one-row-per-group layout.
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = 13
$ws1.Range("C2").Value = "A13"
$ws1.Range("D2").Value = "team1"
$ws1.Range("E2").Value = "rohsurve"

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = 13
$ws1.Range("C3").Value = "A13"
$ws1.Range("D3").Value = "team1"
$ws1.Range("E3").Value = "bsairamr"

$ws1.Range("A4").Value = 1
$ws1.Range("B4").Value = 13
$ws1.Range("C4").Value = "A13"
$ws1.Range("D4").Value = "team1"
$ws1.Range("E4").Value = "yangsis"

$ws1.Range("A5").Value = 1
$ws1.Range("B5").Value = 13
$ws1.Range("C5").Value = "A13"
$ws1.Range("D5").Value = "team1"
$ws1.Range("E5").Value = "aykaranj"
$ws1.Range("F5").Value = 1

$ws1.Range("A6").Value = 1
$ws1.Range("B6").Value = 13
$ws1.Range("C6").Value = "A13"
$ws1.Range("D6").Value = "team1"
$ws1.Range("E6").Value = "voraj"

$ws1.Range("A7").Value = 1
$ws1.Range("B7").Value = 13
$ws1.Range("C7").Value = "A13"
$ws1.Range("D7").Value = "team2"
$ws1.Range("E7").Value = "samvaity"

$ws1.Range("A8").Value = 1
$ws1.Range("B8").Value = 13
$ws1.Range("C8").Value = "A13"
$ws1.Range("D8").Value = "team2"
$ws1.Range("E8").Value = "cdingers"

$ws1.Range("A9").Value = 1
$ws1.Range("B9").Value = 13
$ws1.Range("C9").Value = "A13"
$ws1.Range("D9").Value = "team2"
$ws1.Range("E9").Value = "yketkar"
$ws1.Range("F9").Value = 1

$ws1.Range("A10").Value = 1
$ws1.Range("B10").Value = 13
$ws1.Range("C10").Value = "A13"
$ws1.Range("D10").Value = "team2"
$ws1.Range("E10").Value = "nishshah"

$ws1.Range("A11").Value = 1
$ws1.Range("B11").Value = 13
$ws1.Range("C11").Value = "A13"
$ws1.Range("D11").Value = "team2"
$ws1.Range("E11").Value = "arivero"

$ws1.Range("A12").Value = 1
$ws1.Range("B12").Value = 13
$ws1.Range("C12").Value = "A13"
$ws1.Range("D12").Value = "team3"
$ws1.Range("E12").Value = "sshivara"

$ws1.Range("A13").Value = 1
$ws1.Range("B13").Value = 13
$ws1.Range("C13").Value = "A13"
$ws1.Range("D13").Value = "team3"
$ws1.Range("E13").Value = "pmpande"
$ws1.Range("F13").Value = 1

$ws1.Range("A14").Value = 1
$ws1.Range("B14").Value = 13
$ws1.Range("C14").Value = "A13"
$ws1.Range("D14").Value = "team3"
$ws1.Range("E14").Value = "mmlele"

$ws1.Range("A15").Value = 1
$ws1.Range("B15").Value = 13
$ws1.Range("C15").Value = "A13"
$ws1.Range("D15").Value = "team3"
$ws1.Range("E15").Value = "acsarkis"

$ws1.Range("A16").Value = 1
$ws1.Range("B16").Value = 13
$ws1.Range("C16").Value = "A13"
$ws1.Range("D16").Value = "team4"
$ws1.Range("E16").Value = "vpatani"

$ws1.Range("A17").Value = 1
$ws1.Range("B17").Value = 13
$ws1.Range("C17").Value = "A13"
$ws1.Range("D17").Value = "team4"
$ws1.Range("E17").Value = "harranga"

$ws1.Range("A18").Value = 1
$ws1.Range("B18").Value = 13
$ws1.Range("C18").Value = "A13"
$ws1.Range("D18").Value = "team4"
$ws1.Range("E18").Value = "asadana"
$ws1.Range("F18").Value = 1

$ws1.Range("A19").Value = 1
$ws1.Range("B19").Value = 13
$ws1.Range("C19").Value = "A13"
$ws1.Range("D19").Value = "team4"
$ws1.Range("E19").Value = "jashjhav"

$ws1.Range("A20").Value = 1
$ws1.Range("B20").Value = 13
$ws1.Range("C20").Value = "A13"
$ws1.Range("D20").Value = "team5"
$ws1.Range("E20").Value = "guzh"

$ws1.Range("A21").Value = 1
$ws1.Range("B21").Value = 13
$ws1.Range("C21").Value = "A13"
$ws1.Range("D21").Value = "team5"
$ws1.Range("E21").Value = "deng4"
$ws1.Range("F21").Value = 1

$ws1.Range("A22").Value = 1
$ws1.Range("B22").Value = 13
$ws1.Range("C22").Value = "A13"
$ws1.Range("D22").Value = "team5"
$ws1.Range("E22").Value = "jaynagle"

$ws1.Range("A23").Value = 1
$ws1.Range("B23").Value = 13
$ws1.Range("C23").Value = "A13"
$ws1.Range("D23").Value = "team5"
$ws1.Range("E23").Value = "rnedunur"

# ---------------------------------------------------------------------------
# 3. Best-effort column widths for the two brand new columns (E, F).
# ---------------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 11.02
$ws1.Columns.Item(6).ColumnWidth = 11.88

# ---------------------------------------------------------------------------
# 4. Update the saved selections on both sheets. "group-student" is updated
#    first so that "groups" ends up as the active sheet/tab, matching the
#    original file.
# ---------------------------------------------------------------------------
$ws2.Range("D31").Select() | Out-Null
$ws1.Range("F21").Select() | Out-Null
